$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (55 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 291.2353
$ws.Range("I33").Value = 311.69232
$ws.Range("K33").Value = 311.69232
$ws.Range("M33").Value = -82.69232
$ws.Range("H80").Value = 4810.2856
$ws.Range("I80").Value = 7792.7144
$ws.Range("J80").Value = 1827.8572
$ws.Range("K80").Value = 23378.1432
$ws.Range("L80").Value = 5483.571599999999
$ws.Range("M80").Value = -22380.1432
$ws.Range("N80").Value = -7479.571599999999
$ws.Range("H83").Value = 4810.2856
$ws.Range("I83").Value = 7792.7144
$ws.Range("J83").Value = 1827.8572
$ws.Range("K83").Value = 70134.4296
$ws.Range("L83").Value = 16450.7148
$ws.Range("M83").Value = -65142.4296
$ws.Range("N83").Value = -26434.7148
$ws.Range("H103").Value = 1011.88
$ws.Range("I103").Value = 850.7
$ws.Range("J103").Value = 1119.3334
$ws.Range("K103").Value = 2552.1
$ws.Range("L103").Value = 3358.0002
$ws.Range("M103").Value = -1966.1
$ws.Range("N103").Value = -4530.0002
$ws.Range("H127").Value = 8390.333000000001
$ws.Range("I127").Value = 1841.5714
$ws.Range("J127").Value = 17558.6
$ws.Range("K127").Value = 5524.7142
$ws.Range("L127").Value = 52675.8
$ws.Range("M127").Value = -564.7142000000003
$ws.Range("N127").Value = -62595.8
$ws.Range("H129").Value = 145267.58
$ws.Range("I129").Value = 252159.25
$ws.Range("K129").Value = 756477.75
$ws.Range("M129").Value = -751477.75
$ws.Range("H132").Value = 13542.188
$ws.Range("I132").Value = 3598.558
$ws.Range("K132").Value = 10795.674
$ws.Range("M132").Value = -8265.673999999999
$ws.Range("H135").Value = 9263472
$ws.Range("I135").Value = 11364737
$ws.Range("J135").Value = 17907
$ws.Range("K135").Value = 102282633
$ws.Range("L135").Value = 161163
$ws.Range("M135").Value = -102280098
$ws.Range("N135").Value = -166233
$ws.Range("H138").Value = 3554.5676
$ws.Range("I138").Value = 1862.0625
$ws.Range("K138").Value = 5586.1875
$ws.Range("M138").Value = -446.1875
$ws.Range("H141").Value = 8001
$ws.Range("I141").Value = 4538.375
$ws.Range("K141").Value = 13615.125
$ws.Range("M141").Value = -8435.125

# --- Sheet: ARM (20 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1993.3928
$ws.Range("I61").Value = 1576.125
$ws.Range("K61").Value = 1576.125
$ws.Range("M61").Value = -1364.125
$ws.Range("H74").Value = 2381.6
$ws.Range("I74").Value = 1853.25
$ws.Range("K74").Value = 1853.25
$ws.Range("M74").Value = -979.25
$ws.Range("H77").Value = 2381.6
$ws.Range("I77").Value = 1853.25
$ws.Range("K77").Value = 9266.25
$ws.Range("M77").Value = -4898.25
$ws.Range("H132").Value = 1356.2307
$ws.Range("I132").Value = 1356.2307
$ws.Range("K132").Value = 4068.6921
$ws.Range("M132").Value = -1538.6921
$ws.Range("H136").Value = 1993.3928
$ws.Range("I136").Value = 1576.125
$ws.Range("K136").Value = 4728.375
$ws.Range("M136").Value = -2178.375

# --- Sheet: BSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6411781
$ws.Range("I94").Value = 1795.421
$ws.Range("J94").Value = 23810314
$ws.Range("K94").Value = 1795.421
$ws.Range("L94").Value = 23810314
$ws.Range("M94").Value = -1344.421
$ws.Range("N94").Value = -23811216
$ws.Range("H99").Value = 48724.332
$ws.Range("I99").Value = 61431.285
$ws.Range("K99").Value = 61431.285
$ws.Range("M99").Value = -59933.285
$ws.Range("H107").Value = 3552.9
$ws.Range("I107").Value = 3191.25
$ws.Range("K107").Value = 3191.25
$ws.Range("M107").Value = -1271.25
$ws.Range("H134").Value = 3171.9285
$ws.Range("I134").Value = 2564.5173
$ws.Range("K134").Value = 7693.5519
$ws.Range("M134").Value = -5158.5519

# --- Sheet: CRP (20 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5354.7856
$ws.Range("I16").Value = 6081.636
$ws.Range("K16").Value = 6081.636
$ws.Range("M16").Value = -5794.636
$ws.Range("H31").Value = 1682.1842
$ws.Range("I31").Value = 1373.9667
$ws.Range("K31").Value = 1373.9667
$ws.Range("M31").Value = -1078.9667
$ws.Range("H34").Value = 1682.1842
$ws.Range("I34").Value = 1373.9667
$ws.Range("K34").Value = 1373.9667
$ws.Range("M34").Value = -1171.9667
$ws.Range("H107").Value = 10399.714
$ws.Range("I107").Value = 776
$ws.Range("K107").Value = 776
$ws.Range("M107").Value = 1144
$ws.Range("H113").Value = 5354.7856
$ws.Range("I113").Value = 6081.636
$ws.Range("K113").Value = 6081.636
$ws.Range("M113").Value = -3911.636

# --- Sheet: CUL (35 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1852.75
$ws.Range("I5").Value = 2091.2727
$ws.Range("J5").Value = 1328
$ws.Range("K5").Value = 6273.8181
$ws.Range("L5").Value = 3984
$ws.Range("M5").Value = -6161.8181
$ws.Range("N5").Value = -4208
$ws.Range("H107").Value = 332.2
$ws.Range("I107").Value = 345
$ws.Range("J107").Value = 323.66666
$ws.Range("K107").Value = 1035
$ws.Range("L107").Value = 970.9999799999999
$ws.Range("M107").Value = 885
$ws.Range("N107").Value = -4810.99998
$ws.Range("H121").Value = 1386.7894
$ws.Range("I121").Value = 291
$ws.Range("J121").Value = 2604.3333
$ws.Range("K121").Value = 873
$ws.Range("L121").Value = 7812.999899999999
$ws.Range("M121").Value = 437
$ws.Range("N121").Value = -10432.9999
$ws.Range("H132").Value = 1138.2222
$ws.Range("I132").Value = 1359
$ws.Range("J132").Value = 862.25
$ws.Range("K132").Value = 12231
$ws.Range("L132").Value = 7760.25
$ws.Range("M132").Value = -9701
$ws.Range("N132").Value = -12820.25
$ws.Range("H135").Value = 1852.75
$ws.Range("I135").Value = 2091.2727
$ws.Range("J135").Value = 1328
$ws.Range("K135").Value = 18821.4543
$ws.Range("L135").Value = 11952
$ws.Range("M135").Value = -16286.4543
$ws.Range("N135").Value = -17022

# --- Sheet: GSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 34999.5
$ws.Range("J52").Value = 34999.5
$ws.Range("L52").Value = 34999.5
$ws.Range("N52").Value = -35517.5
$ws.Range("H80").Value = 23891522
$ws.Range("J80").Value = 30305974
$ws.Range("L80").Value = 30305974
$ws.Range("N80").Value = -30307970
$ws.Range("H83").Value = 23891522
$ws.Range("J83").Value = 30305974
$ws.Range("L83").Value = 151529870
$ws.Range("N83").Value = -151539854

# --- Sheet: LTW (19 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 44493.332
$ws.Range("J57").Value = 44493.332
$ws.Range("L57").Value = 44493.332
$ws.Range("N57").Value = -45625.332
$ws.Range("H117").Value = 51195.5
$ws.Range("J117").Value = 51195.5
$ws.Range("L117").Value = 51195.5
$ws.Range("N117").Value = -60373.5
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
$ws.Range("H136").Value = 2455.4893
$ws.Range("I136").Value = 2284.3865
$ws.Range("J136").Value = 4965
$ws.Range("K136").Value = 6853.1595
$ws.Range("L136").Value = 14895
$ws.Range("M136").Value = -4303.1595
$ws.Range("N136").Value = -19995

# --- Sheet: WVR (11 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("M18").Value = -9827
$ws.Range("H126").Value = 2791.6875
$ws.Range("J126").Value = 3301.8
$ws.Range("L126").Value = 9905.400000000001
$ws.Range("N126").Value = -14845.4
